$wb = $excel.ActiveWorkbook

$wsOpt = $wb.Worksheets.Item("optimization_parameters")

# Row 1: drop the trailing duplicate "value" header cells (C1:F1), keep A1/B1.
$wsOpt.Range("C1:F1").ClearContents()

# Row 8: "Model" -> "production_function" (value stays "Sigmoid").
$wsOpt.Range("A8").Value = "production_function"

# Insert a new row 9 for the "L_curve" parameter (pushes the rest down by one).
$wsOpt.Rows.Item(9).Insert()
$wsOpt.Range("A9").Value = "L_curve"
$wsOpt.Range("B9").Value = 1
$wsOpt.Range("B9").NumberFormat = "0.00E+00"

# The old "Deletion" row (now shifted down to row 17) is removed entirely.
$wsOpt.Rows.Item(17).Delete()

# Make optimization_parameters the active sheet / tab, with the new selection.
$wsOpt.Activate()
$wsOpt.Range("C20").Select()
